# Generate Report for Handoff
#
# Refreshes the handoff timestamps for the file
# "2e866cb6-0967-44c0-a16d-58105418b33b.md" (table row 5 on every sheet),
# which is in "Ready for handoff" status, to reflect a newly generated
# handoff report.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest Handoff Date" column (D) for row 5.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("D5").Value = "2016-03-23 04:43:30"

# zh-cn sheet: "Latest Handoff Datetime" column (E) for row 5.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E5").Value = "2016-03-23 04:43:26"

# de-de sheet: "Latest Handoff Datetime" column (E) for row 5.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E5").Value = "2016-03-23 04:43:30"
